$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: the paragraph that starts "-El boton de pagina principal
# nunca da los eventos..." gets a one-letter spelling fix (util ->
# util with an accented u) and the run carrying the text is split
# into three runs at the edit boundaries, the way Word leaves a run
# split behind at the points where text was retyped.
# -----------------------------------------------------------------
$p4 = $d.Paragraphs(4).Range
$start4 = $p4.Start
$full4 = $p4.Text

# Fix the accent: "util" -> "util" with u-acute. Anchored on the
# ASCII-only "algo u" so the literal accented characters elsewhere in
# the paragraph don't need to be retyped here.
$anchorFix = "algo u"
$uIdx = $full4.IndexOf($anchorFix) + $anchorFix.Length - 1
$fixRange = $d.Range($start4 + $uIdx, $start4 + $uIdx + 1)
$fixRange.Text = [string][char]0xFA

# Split positions, all anchored on ASCII-only substrings so they are
# not sensitive to how this file's own encoding round-trips:
#   after "nunca da"                -> "...nunca da" | " los eventos..."
#   after "algo u"                  -> "...algo u"    | "til antes..."
#   before ". Simile con gente"     -> "...pestanas)" | ". Simile..."
#   after ". Simile con gente"      -> "...gente"     | "."
$daIdx = $full4.IndexOf("nunca da") + "nunca da".Length
$tilIdx = $full4.IndexOf("algo u") + "algo u".Length
$simileAnchor = ". Simile con gente"
$pestIdx = $full4.IndexOf($simileAnchor)
$genteIdx = $pestIdx + $simileAnchor.Length

# Apply the splits right-to-left so earlier offsets stay valid; each
# split is produced by dropping a temporary bookmark at the boundary
# (OOXML requires bookmarks to sit between runs) and then removing it
# again, which leaves the run break behind without leaving a bookmark.
foreach ($splitAt in @($genteIdx, $pestIdx, $tilIdx, $daIdx)) {
    $boundary = $d.Range($start4 + $splitAt, $start4 + $splitAt)
    $d.Bookmarks.Add("ZZZTempSplit", $boundary)
    $d.Bookmarks("ZZZTempSplit").Delete()
}

# -----------------------------------------------------------------
# Change 2: the paragraph "Es posible que sea algo como lo q pasaba
# en servicios..." has its run split right after "estruct" / before
# "ura" -- this is exactly where the document's _GoBack bookmark
# (marking the last edit position) needs to move to.
# -----------------------------------------------------------------
$p5 = $d.Paragraphs(5).Range
$start5 = $p5.Start
$full5 = $p5.Text

$splitIdx = $full5.IndexOf("estructura") + "estruct".Length

$d.Bookmarks("_GoBack").Delete()
$goBack = $d.Range($start5 + $splitIdx, $start5 + $splitIdx)
$d.Bookmarks.Add("_GoBack", $goBack)
